# Add a new "2021" data column (column O) to the right of the existing
# 2010-2020 table, mirroring the per-row formatting already used in
# column N, then update the sheet view's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the column (rather than copy/paste-special into a bare range)
# makes Excel clone each row's existing N-column style onto the new O
# column cell-for-cell, which keeps the exact style actually used instead
# of silently collapsing to a different-but-visually-equal style id.
$ws.Columns("O:O").Insert()

# Row 15 is a spacer/section row in the source table and never receives a
# 2021 figure, so drop the blank placeholder cell the column insert left
# behind there.
$ws.Range("O15").Clear()

# Row -> 2021 value, keyed by the cell address that receives it.
$values = [ordered]@{
    "O4"  = 2021
    "O5"  = 11.7
    "O6"  = 16.4
    "O7"  = 9.7
    "O8"  = 12.1
    "O9"  = 5.3
    "O10" = 4.7
    "O11" = 3.4
    "O12" = 18.8
    "O13" = 19.6
    "O14" = 6.9
    "O16" = 12.8
    "O17" = 11
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Move the active selection/view to R11 (this also clears the old
# topLeftCell scroll-freeze that pointed at E2).
$ws.Range("R11").Select()
